# Unify the conception of DataNode, DataTable, Entity:
# rename the two worksheets and make the second (DataTable) the active tab,
# matching the author's re-save where "DataTable" ends up selected.

$wb = $excel.ActiveWorkbook

$wsProperty = $wb.Worksheets.Item("Property1")
$wsRecord   = $wb.Worksheets.Item("Record")

$wsProperty.Name = "DataNode"
$wsRecord.Name   = "DataTable"

$wsRecord.Activate()
